# Updates "posicion_jugadores" sheet data through Jornada 4 (up to J5 per commit message).
# The player roster rows (2-8) shifted down by one position (a new player - Franco
# Zanelatto - now leads the list, and each previously listed player moved down one
# row), and several userCount (column F) values were refreshed with newer counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Franco Zanelatto (firstName/lastName stay blank) ---
$ws.Range("A2").Value = "Franco Zanelatto"
$ws.Range("B2").Value = "franco-zanelatto"
$ws.Range("C2").Value = "F. Zanelatto"
$ws.Range("D2").Value = "M"
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 330
$ws.Range("G2").Value = 973650
$ws.Range("J2").Value = 86.433333333333
$ws.Range("K2").Value = 79.666666666667
$ws.Range("L2").Value = 3

# --- Row 3: Ricardo Lagos (firstName/lastName become blank) ---
$ws.Range("A3").Value = "Ricardo Lagos"
$ws.Range("B3").Value = "ricardo-lagos"
$ws.Range("C3").Value = "R. Lagos"
$ws.Range("D3").Value = "D"
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = 92
$ws.Range("G3").Value = 973682
$ws.Range("H3").ClearContents()
$ws.Range("J3").Value = 57.991836734694
$ws.Range("K3").Value = 81.477551020408
$ws.Range("L3").Value = 49

# --- Row 4: Ángel De la Cruz (firstName keeps a value, now his own name) ---
$ws.Range("A4").Value = "Ángel De la Cruz"
$ws.Range("B4").Value = "angel-de-la-cruz"
$ws.Range("C4").Value = "Á. D. l. Cruz"
$ws.Range("D4").Value = "G"
$ws.Range("E4").Value = 12
$ws.Range("F4").Value = 13
$ws.Range("G4").Value = 1109887
$ws.Range("H4").Value = "Ángel De la Cruz"
$ws.Range("J4").Value = 8.921739130434799
$ws.Range("K4").Value = 51.747826086957
$ws.Range("L4").Value = 46

# --- Row 5: Catriel Cabellos (firstName now populated with his own name) ---
$ws.Range("A5").Value = "Catriel Cabellos"
$ws.Range("B5").Value = "cabellos-catriel"
$ws.Range("C5").Value = "C. Cabellos"
$ws.Range("D5").Value = "M"
$ws.Range("E5").Value = 27
$ws.Range("F5").Value = 320
$ws.Range("G5").Value = 1415942
$ws.Range("H5").Value = "Catriel Cabellos"
$ws.Range("J5").Value = 45.54358974359
$ws.Range("K5").Value = 19.328205128205
$ws.Range("L5").Value = 39

# --- Row 6: Jesús Castillo (firstName/lastName stay blank) ---
$ws.Range("A6").Value = "Jesús Castillo"
$ws.Range("B6").Value = "jesus-castillo"
$ws.Range("C6").Value = "J. Castillo"
$ws.Range("D6").Value = "M"
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = 81
$ws.Range("G6").Value = 913398
$ws.Range("J6").Value = 51.546153846154
$ws.Range("K6").Value = 39.507692307692
$ws.Range("L6").Value = 13

# --- Row 7: Juan Freytes (firstName/lastName stay blank) ---
$ws.Range("A7").Value = "Juan Freytes"
$ws.Range("B7").Value = "juan-freytes"
$ws.Range("C7").Value = "J. Freytes"
$ws.Range("D7").Value = "D"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 49
$ws.Range("G7").Value = 962187
$ws.Range("J7").Value = 39.205172413793
$ws.Range("K7").Value = 64.524137931035
$ws.Range("L7").Value = 58

# --- Row 8: Aldair Fuentes (firstName/lastName stay blank) ---
$ws.Range("A8").Value = "Aldair Fuentes"
$ws.Range("B8").Value = "aldair-fuentes"
$ws.Range("C8").Value = "A. Fuentes"
$ws.Range("D8").Value = "M"
$ws.Range("E8").Value = 20
$ws.Range("F8").Value = 132
$ws.Range("G8").Value = 876927
$ws.Range("J8").Value = 29.65
$ws.Range("K8").Value = 88.65000000000001
$ws.Range("L8").Value = 2

# --- userCount refreshes for remaining roster rows ---
$ws.Range("F9").Value = 215
$ws.Range("F11").Value = 468
$ws.Range("F12").Value = 266
$ws.Range("F13").Value = 103
$ws.Range("F14").Value = 72
$ws.Range("F15").Value = 84
$ws.Range("F16").Value = 585
